$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell references to their new values, derived from the commit diff.
$updates = [ordered]@{
    'D2' = '30.499.15'
    'E2' = '  -1.05%  '
    'D3' = '2.112.54'
    'E3' = '  -0.23%  '
    'E4' = '  +0.12%  '
    'D5' = '335.06'
    'E6' = '  +0.14%  '
    'D7' = '0.5254'
    'E7' = '  -1.46%  '
    'D8' = '0.4495'
    'E8' = '  +2.20%  '
    'D9' = '53.41'
    'E9' = '  +12.67%  '
    'D10' = '0.09023'
    'E10' = '  +0.07%  '
    'D11' = '1.174'
    'E11' = '  -0.79%  '
    'D12' = '24.52'
    'E12' = '  -1.95%  '
    'D13' = '2.102.99'
    'E13' = '  -0.39%  '
    'D14' = '6.790'
    'E14' = '  +0.24%  '
    'D15' = '7.826'
    'E15' = '  +0.03%  '
    'D16' = '96.87'
    'E16' = '  +0.03%  '
    'D18' = '0.00001129'
    'E18' = '  -0.38%  '
    'D19' = '0.06620'
    'E19' = '  -0.89%  '
    'D20' = '19.47'
    'E20' = '  +1.60%  '
    'E21' = '  +0.12%  '
    'E22' = '  -0.47%  '
    'D23' = '30.543.85'
    'E23' = '  -1.10%  '
    'E24' = '  +0.59%  '
    'D25' = '2.355'
    'E25' = '  +2.70%  '
    'D26' = '2.353.30'
    'E26' = '  -0.38%  '
    'D27' = '22.44'
    'E27' = '  -1.62%  '
    'D28' = '2.595'
    'E28' = '  -0.12%  '
    'D29' = '163.62'
    'E29' = '  +0.18%  '
    'D30' = '132.89'
    'E30' = '  -0.43%  '
    'D31' = '1.202'
    'E31' = '  +1.76%  '
    'D32' = '0.1076'
    'E32' = '  -0.61%  '
    'D33' = '1.668'
    'E33' = '  +7.67%  '
    'D34' = '6.175'
    'E34' = '  -1.09%  '
    'D35' = '3.935'
    'E35' = '  -1.96%  '
    'D36' = '10.63'
    'E36' = '  +12.09%  '
    'D37' = '0.02581'
    'E37' = '  -0.84%  '
    'D38' = '0.06847'
    'E38' = '  +1.22%  '
    'D39' = '5.567'
    'E39' = '  +0.16%  '
    'D40' = '12.82'
    'E40' = '  -0.77%  '
    'D41' = '0.2302'
    'E41' = '  +0.41%  '
    'D42' = '0.6953'
    'E42' = '  +1.46%  '
    'D43' = '1.247'
    'E43' = '  -0.18%  '
    'D44' = '2.393'
    'E44' = '  +7.24%  '
    'D45' = '1.001'
    'E45' = '  +0.13%  '
    'B46' = 'EnergySwap'
    'C46' = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
    'D46' = '14.11'
    'E46' = '  -0.48%  '
    'B47' = 'Decentraland'
    'C47' = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
    'D47' = '0.6413'
    'E47' = '  -0.67%  '
    'D48' = '3.669'
    'E48' = '  +0.15%  '
    'D49' = '1.252'
    'E49' = '  -2.05%  '
    'D50' = '1.214'
    'E50' = '  +3.43%  '
    'D51' = '83.44'
    'E51' = '  +0.25%  '
}

foreach ($cellRef in $updates.Keys) {
    $cell = $ws.Range($cellRef)
    # Force text format so numeric-looking strings (prices, percents) keep
    # their exact original formatting (leading zeros, trailing zeros, dot
    # groupings, spacing, percent signs) instead of being parsed as numbers.
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$cellRef]
}

Write-Host "Applied $($updates.Count) cell updates"
